# Add a new "Transcriptomics" tag (column F) to the Tags block of the
# minSCe isa_template sheet (rows 13-17), mirroring the existing OBI-backed
# tags already present in columns B/C/D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

# New tag term
$ws.Range("F13").Value = "Transcriptomics"
# New tag term accession number
$ws.Range("F14").Value = "https://bioregistry.io/NCIT:C153189"
# New tag term source REF
$ws.Range("F15").Value = "NCIT"
# New tag description (Comment[description])
$ws.Range("F16").Value = '"A study of the complete set of RNA transcripts that are produced by the genome, under specific circumstances or in a specific cell." []'
# New tag Comment[isObsolete] (force text, not boolean, to match sibling cells)
$ws.Range("F17").Value = "'false"

# Match the row height used by the rest of the Tags rows
$ws.Range("A13:F17").RowHeight = 13.8

# Update the active selection to the newly added column, like the source edit
$ws.Activate()
$ws.Range("F13:F17").Select()
